{"js": "// Replace the math-problem answers in the single table of the document.\n// Mapping extracted from the authoritative diff: [row, col, oldText, newText]\nconst replacements = [\n  [0, 0, \"45-33=12\", \"23+62=85\"],\n  [0, 1, \"29+6=35\", \"54-50=4\"],\n  [0, 2, \"33+6=39\", \"57+25=82\"],\n  [0, 3, \"92-59=33\", \"64+7=71\"],\n  [0, 4, \"99-84=15\", \"12-3=9\"],\n  [1, 0, \"13-11=2\", \"88+8=96\"],\n  [1, 1, \"25+43=68\", \"15-3=12\"],\n  [1, 2, \"20+34=54\", \"41+4=45\"],\n  [1, 3, \"20-6=14\", \"67-61=6\"],\n  [1, 4, \"70-64=6\", \"51-25=26\"],\n  [2, 0, \"74-39=35\", \"51-0=51\"],\n  [2, 1, \"36-24=12\", \"38-12=26\"],\n  [2, 2, \"16+25=41\", \"39+15=54\"],\n  [2, 3, \"72+9=81\", \"9+34=43\"],\n  [2, 4, \"97-78=19\", \"69+28=97\"],\n  [3, 0, \"44+17=61\", \"37+19=56\"],\n  [3, 1, \"70+27=97\", \"8+65=73\"],\n  [3, 2, \"60-21=39\", \"16-15=1\"],\n  [3, 3, \"57+30=87\", \"0+21=21\"],\n  [3, 4, \"17+29=46\", \"85-45=40\"],\n  [4, 0, \"18+1=19\", \"5+76=81\"],\n  [4, 1, \"28+53=81\", \"87-24=63\"],\n  [4, 2, \"38+15=53\", \"29-22=7\"],\n  [4, 3, \"49+42=91\", \"4+85=89\"],\n  [4, 4, \"8+26=34\", \"7+82=89\"],\n  [5, 0, \"0+88=88\", \"15+0=15\"],\n  [5, 1, \"63-61=2\", \"66-35=31\"],\n  [5, 2, \"19-4=15\", \"60-27=33\"],\n  [5, 3, \"5+27=32\", \"49+35=84\"],\n  [5, 4, \"94-52=42\", \"93-8=85\"],\n  [6, 0, \"38+54=92\", \"27-6=21\"],\n  [6, 1, \"20+9=29\", \"80-60=20\"],\n  [6, 2, \"19-18=1\", \"40+47=87\"],\n  [6, 3, \"18-13=5\", \"47+40=87\"],\n  [6, 4, \"14+44=58\", \"88+11=99\"],\n  [7, 0, \"23+17=40\", \"25+57=82\"],\n  [7, 1, \"18+17=35\", \"13+57=70\"],\n  [7, 2, \"27-19=8\", \"17-2=15\"],\n  [7, 3, \"81-27=54\", \"85+9=94\"],\n  [7, 4, \"4+47=51\", \"79-44=35\"],\n  [8, 0, \"17+36=53\", \"38+28=66\"],\n  [8, 1, \"76-7=69\", \"42-23=19\"],\n  [8, 2, \"16+66=82\", \"66-47=19\"],\n  [8, 3, \"35+31=66\", \"12+48=60\"],\n  [8, 4, \"81-37=44\", \"70-3=67\"],\n  [9, 0, \"82+2=84\", \"25+7=32\"],\n  [9, 1, \"53-48=5\", \"77-54=23\"],\n  [9, 2, \"25+14=39\", \"93-20=73\"],\n  [9, 3, \"18-0=18\", \"63-30=33\"],\n  [9, 4, \"42-3=39\", \"6+28=34\"],\n  [10, 0, \"34+27=61\", \"65-18=47\"],\n  [10, 1, \"48-25=23\", \"21+47=68\"],\n  [10, 2, \"95-49=46\", \"32-21=11\"],\n  [10, 3, \"73+11=84\", \"91-41=50\"],\n  [10, 4, \"55-27=28\", \"47+7=54\"],\n  [11, 0, \"32-16=16\", \"91-36=55\"],\n  [11, 1, \"96+0=96\", \"20+33=53\"],\n  [11, 2, \"83-37=46\", \"63-27=36\"],\n  [11, 3, \"54-25=29\", \"4+4=8\"],\n  [11, 4, \"91+0=91\", \"93-28=65\"],\n  [12, 0, \"85-79=6\", \"31+40=71\"],\n  [12, 1, \"96-42=54\", \"28-3=25\"],\n  [12, 2, \"40+28=68\", \"60-59=1\"],\n  [12, 3, \"41+31=72\", \"79-70=9\"],\n  [12, 4, \"27-6=21\", \"54+41=95\"],\n  [13, 0, \"27-11=16\", \"86+0=86\"],\n  [13, 1, \"65+19=84\", \"15-0=15\"],\n  [13, 2, \"26+35=61\", \"65+31=96\"],\n  [13, 3, \"8+44=52\", \"40+19=59\"],\n  [13, 4, \"91-60=31\", \"12-9=3\"],\n  [14, 0, \"75-43=32\", \"3+28=31\"],\n  [14, 1, \"93-10=83\", \"22-20=2\"],\n  [14, 2, \"79-43=36\", \"40+8=48\"],\n  [14, 3, \"26-23=3\", \"34+56=90\"],\n  [14, 4, \"91-57=34\", \"41-4=37\"],\n  [15, 0, \"90-89=1\", \"48-37=11\"],\n  [15, 1, \"28-11=17\", \"28+49=77\"],\n  [15, 2, \"44-18=26\", \"14+75=89\"],\n  [15, 3, \"12-10=2\", \"7+55=62\"],\n  [15, 4, \"57+3=60\", \"75-51=24\"],\n  [16, 0, \"25-15=10\", \"53+39=92\"],\n  [16, 1, \"95-20=75\", \"59-30=29\"],\n  [16, 2, \"3+86=89\", \"98-2=96\"],\n  [16, 3, \"26+46=72\", \"48+16=64\"],\n  [16, 4, \"20+35=55\", \"41+32=73\"],\n  [17, 0, \"54-6=48\", \"41-8=33\"],\n  [17, 1, \"10+46=56\", \"74-8=66\"],\n  [17, 2, \"4+40=44\", \"34+52=86\"],\n  [17, 3, \"18+7=25\", \"5+2=7\"],\n  [17, 4, \"97-0=97\", \"96-24=72\"],\n  [18, 0, \"85-30=55\", \"7+87=94\"],\n  [18, 1, \"90-67=23\", \"23+28=51\"],\n  [18, 2, \"58+7=65\", \"57+22=79\"],\n  [18, 3, \"87-16=71\", \"47-29=18\"],\n  [18, 4, \"66-14=52\", \"82-75=7\"],\n  [19, 0, \"84-59=25\", \"19+26=45\"],\n  [19, 1, \"15+71=86\", \"21+68=89\"],\n  [19, 2, \"51-36=15\", \"36+37=73\"],\n  [19, 3, \"70+25=95\", \"11+30=41\"],\n  [19, 4, \"22+42=64\", \"39-25=14\"],\n];\n\nconst table = context.document.body.tables.getFirst();\nawait context.sync();\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange();\n  range.load(\"text\");\n  await context.sync();\n  // Word represents a table-cell range's text with a trailing cell-mark\n  // character (tab), so compare against the trimmed text.\n  const currentText = range.text.replace(/\\t$/, \"\");\n  if (currentText !== oldText) {\n    throw new Error(`Unexpected cell text at (${row},${col}): got \"${currentText}\", expected \"${oldText}\"`);\n  }\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the math-problem answers in the single table of the document.\n# Mapping extracted from the authoritative diff (1-based Row/Col indices\n# matching Word COM's Table.Cell(row, col) convention).\n$replacements = @(\n    @{Row=1; Col=1; Old=\"45-33=12\"; New=\"23+62=85\"},\n    @{Row=1; Col=2; Old=\"29+6=35\"; New=\"54-50=4\"},\n    @{Row=1; Col=3; Old=\"33+6=39\"; New=\"57+25=82\"},\n    @{Row=1; Col=4; Old=\"92-59=33\"; New=\"64+7=71\"},\n    @{Row=1; Col=5; Old=\"99-84=15\"; New=\"12-3=9\"},\n    @{Row=2; Col=1; Old=\"13-11=2\"; New=\"88+8=96\"},\n    @{Row=2; Col=2; Old=\"25+43=68\"; New=\"15-3=12\"},\n    @{Row=2; Col=3; Old=\"20+34=54\"; New=\"41+4=45\"},\n    @{Row=2; Col=4; Old=\"20-6=14\"; New=\"67-61=6\"},\n    @{Row=2; Col=5; Old=\"70-64=6\"; New=\"51-25=26\"},\n    @{Row=3; Col=1; Old=\"74-39=35\"; New=\"51-0=51\"},\n    @{Row=3; Col=2; Old=\"36-24=12\"; New=\"38-12=26\"},\n    @{Row=3; Col=3; Old=\"16+25=41\"; New=\"39+15=54\"},\n    @{Row=3; Col=4; Old=\"72+9=81\"; New=\"9+34=43\"},\n    @{Row=3; Col=5; Old=\"97-78=19\"; New=\"69+28=97\"},\n    @{Row=4; Col=1; Old=\"44+17=61\"; New=\"37+19=56\"},\n    @{Row=4; Col=2; Old=\"70+27=97\"; New=\"8+65=73\"},\n    @{Row=4; Col=3; Old=\"60-21=39\"; New=\"16-15=1\"},\n    @{Row=4; Col=4; Old=\"57+30=87\"; New=\"0+21=21\"},\n    @{Row=4; Col=5; Old=\"17+29=46\"; New=\"85-45=40\"},\n    @{Row=5; Col=1; Old=\"18+1=19\"; New=\"5+76=81\"},\n    @{Row=5; Col=2; Old=\"28+53=81\"; New=\"87-24=63\"},\n    @{Row=5; Col=3; Old=\"38+15=53\"; New=\"29-22=7\"},\n    @{Row=5; Col=4; Old=\"49+42=91\"; New=\"4+85=89\"},\n    @{Row=5; Col=5; Old=\"8+26=34\"; New=\"7+82=89\"},\n    @{Row=6; Col=1; Old=\"0+88=88\"; New=\"15+0=15\"},\n    @{Row=6; Col=2; Old=\"63-61=2\"; New=\"66-35=31\"},\n    @{Row=6; Col=3; Old=\"19-4=15\"; New=\"60-27=33\"},\n    @{Row=6; Col=4; Old=\"5+27=32\"; New=\"49+35=84\"},\n    @{Row=6; Col=5; Old=\"94-52=42\"; New=\"93-8=85\"},\n    @{Row=7; Col=1; Old=\"38+54=92\"; New=\"27-6=21\"},\n    @{Row=7; Col=2; Old=\"20+9=29\"; New=\"80-60=20\"},\n    @{Row=7; Col=3; Old=\"19-18=1\"; New=\"40+47=87\"},\n    @{Row=7; Col=4; Old=\"18-13=5\"; New=\"47+40=87\"},\n    @{Row=7; Col=5; Old=\"14+44=58\"; New=\"88+11=99\"},\n    @{Row=8; Col=1; Old=\"23+17=40\"; New=\"25+57=82\"},\n    @{Row=8; Col=2; Old=\"18+17=35\"; New=\"13+57=70\"},\n    @{Row=8; Col=3; Old=\"27-19=8\"; New=\"17-2=15\"},\n    @{Row=8; Col=4; Old=\"81-27=54\"; New=\"85+9=94\"},\n    @{Row=8; Col=5; Old=\"4+47=51\"; New=\"79-44=35\"},\n    @{Row=9; Col=1; Old=\"17+36=53\"; New=\"38+28=66\"},\n    @{Row=9; Col=2; Old=\"76-7=69\"; New=\"42-23=19\"},\n    @{Row=9; Col=3; Old=\"16+66=82\"; New=\"66-47=19\"},\n    @{Row=9; Col=4; Old=\"35+31=66\"; New=\"12+48=60\"},\n    @{Row=9; Col=5; Old=\"81-37=44\"; New=\"70-3=67\"},\n    @{Row=10; Col=1; Old=\"82+2=84\"; New=\"25+7=32\"},\n    @{Row=10; Col=2; Old=\"53-48=5\"; New=\"77-54=23\"},\n    @{Row=10; Col=3; Old=\"25+14=39\"; New=\"93-20=73\"},\n    @{Row=10; Col=4; Old=\"18-0=18\"; New=\"63-30=33\"},\n    @{Row=10; Col=5; Old=\"42-3=39\"; New=\"6+28=34\"},\n    @{Row=11; Col=1; Old=\"34+27=61\"; New=\"65-18=47\"},\n    @{Row=11; Col=2; Old=\"48-25=23\"; New=\"21+47=68\"},\n    @{Row=11; Col=3; Old=\"95-49=46\"; New=\"32-21=11\"},\n    @{Row=11; Col=4; Old=\"73+11=84\"; New=\"91-41=50\"},\n    @{Row=11; Col=5; Old=\"55-27=28\"; New=\"47+7=54\"},\n    @{Row=12; Col=1; Old=\"32-16=16\"; New=\"91-36=55\"},\n    @{Row=12; Col=2; Old=\"96+0=96\"; New=\"20+33=53\"},\n    @{Row=12; Col=3; Old=\"83-37=46\"; New=\"63-27=36\"},\n    @{Row=12; Col=4; Old=\"54-25=29\"; New=\"4+4=8\"},\n    @{Row=12; Col=5; Old=\"91+0=91\"; New=\"93-28=65\"},\n    @{Row=13; Col=1; Old=\"85-79=6\"; New=\"31+40=71\"},\n    @{Row=13; Col=2; Old=\"96-42=54\"; New=\"28-3=25\"},\n    @{Row=13; Col=3; Old=\"40+28=68\"; New=\"60-59=1\"},\n    @{Row=13; Col=4; Old=\"41+31=72\"; New=\"79-70=9\"},\n    @{Row=13; Col=5; Old=\"27-6=21\"; New=\"54+41=95\"},\n    @{Row=14; Col=1; Old=\"27-11=16\"; New=\"86+0=86\"},\n    @{Row=14; Col=2; Old=\"65+19=84\"; New=\"15-0=15\"},\n    @{Row=14; Col=3; Old=\"26+35=61\"; New=\"65+31=96\"},\n    @{Row=14; Col=4; Old=\"8+44=52\"; New=\"40+19=59\"},\n    @{Row=14; Col=5; Old=\"91-60=31\"; New=\"12-9=3\"},\n    @{Row=15; Col=1; Old=\"75-43=32\"; New=\"3+28=31\"},\n    @{Row=15; Col=2; Old=\"93-10=83\"; New=\"22-20=2\"},\n    @{Row=15; Col=3; Old=\"79-43=36\"; New=\"40+8=48\"},\n    @{Row=15; Col=4; Old=\"26-23=3\"; New=\"34+56=90\"},\n    @{Row=15; Col=5; Old=\"91-57=34\"; New=\"41-4=37\"},\n    @{Row=16; Col=1; Old=\"90-89=1\"; New=\"48-37=11\"},\n    @{Row=16; Col=2; Old=\"28-11=17\"; New=\"28+49=77\"},\n    @{Row=16; Col=3; Old=\"44-18=26\"; New=\"14+75=89\"},\n    @{Row=16; Col=4; Old=\"12-10=2\"; New=\"7+55=62\"},\n    @{Row=16; Col=5; Old=\"57+3=60\"; New=\"75-51=24\"},\n    @{Row=17; Col=1; Old=\"25-15=10\"; New=\"53+39=92\"},\n    @{Row=17; Col=2; Old=\"95-20=75\"; New=\"59-30=29\"},\n    @{Row=17; Col=3; Old=\"3+86=89\"; New=\"98-2=96\"},\n    @{Row=17; Col=4; Old=\"26+46=72\"; New=\"48+16=64\"},\n    @{Row=17; Col=5; Old=\"20+35=55\"; New=\"41+32=73\"},\n    @{Row=18; Col=1; Old=\"54-6=48\"; New=\"41-8=33\"},\n    @{Row=18; Col=2; Old=\"10+46=56\"; New=\"74-8=66\"},\n    @{Row=18; Col=3; Old=\"4+40=44\"; New=\"34+52=86\"},\n    @{Row=18; Col=4; Old=\"18+7=25\"; New=\"5+2=7\"},\n    @{Row=18; Col=5; Old=\"97-0=97\"; New=\"96-24=72\"},\n    @{Row=19; Col=1; Old=\"85-30=55\"; New=\"7+87=94\"},\n    @{Row=19; Col=2; Old=\"90-67=23\"; New=\"23+28=51\"},\n    @{Row=19; Col=3; Old=\"58+7=65\"; New=\"57+22=79\"},\n    @{Row=19; Col=4; Old=\"87-16=71\"; New=\"47-29=18\"},\n    @{Row=19; Col=5; Old=\"66-14=52\"; New=\"82-75=7\"},\n    @{Row=20; Col=1; Old=\"84-59=25\"; New=\"19+26=45\"},\n    @{Row=20; Col=2; Old=\"15+71=86\"; New=\"21+68=89\"},\n    @{Row=20; Col=3; Old=\"51-36=15\"; New=\"36+37=73\"},\n    @{Row=20; Col=4; Old=\"70+25=95\"; New=\"11+30=41\"},\n    @{Row=20; Col=5; Old=\"22+42=64\"; New=\"39-25=14\"},\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nforeach ($item in $replacements) {\n    $cell = $tbl.Cell($item.Row, $item.Col)\n    $r = $cell.Range\n    # Word COM reports a table cell's Range.Text with the trailing\n    # cell-mark characters (CR + BEL), so trim before comparing.\n    $current = $r.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $item.Old) {\n        throw \"Unexpected cell text at ($($item.Row),$($item.Col)): got '$current', expected '$($item.Old)'\"\n    }\n    $r.Text = $item.New\n}\n"}
